$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rebuild rows 2-31: columns A (timestamp) and B (label) are row-indexed and unchanged;
# columns C-H (ax,ay,az,gx,gy,gz) are updated to their new sensor readings, including the
# newly appended rows 22-31 (timestamps 2000-2900) captured on May 9th.

# row 2 (timestamp 0)
$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 2).Value = "walkingToRunning"
$ws.Cells.Item(2, 3).Value = 7.345117568969727
$ws.Cells.Item(2, 4).Value = -12.58289909362793
$ws.Cells.Item(2, 5).Value = 3.90805721282959
$ws.Cells.Item(2, 6).Value = -0.06342706156940325
$ws.Cells.Item(2, 7).Value = 0.2977593003249712
$ws.Cells.Item(2, 8).Value = 0.5447398664020925

# row 3 (timestamp 100)
$ws.Cells.Item(3, 1).Value = 100
$ws.Cells.Item(3, 2).Value = "walkingToRunning"
$ws.Cells.Item(3, 3).Value = 10.5408239364624
$ws.Cells.Item(3, 4).Value = 2.647232055664062
$ws.Cells.Item(3, 5).Value = 14.10661315917969
$ws.Cells.Item(3, 6).Value = -1.19702434539795
$ws.Cells.Item(3, 7).Value = 1.756288290023804
$ws.Cells.Item(3, 8).Value = 2.010953009128571

# row 4 (timestamp 200)
$ws.Cells.Item(4, 1).Value = 200
$ws.Cells.Item(4, 2).Value = "walkingToRunning"
$ws.Cells.Item(4, 3).Value = -6.684710502624512
$ws.Cells.Item(4, 4).Value = -10.58095169067383
$ws.Cells.Item(4, 5).Value = -20.94781112670898
$ws.Cells.Item(4, 6).Value = 1.257760638143962
$ws.Cells.Item(4, 7).Value = -0.1820772450144785
$ws.Cells.Item(4, 8).Value = -2.106898115902421

# row 5 (timestamp 300)
$ws.Cells.Item(5, 1).Value = 300
$ws.Cells.Item(5, 2).Value = "walkingToRunning"
$ws.Cells.Item(5, 3).Value = -9.718008995056152
$ws.Cells.Item(5, 4).Value = -12.06524276733398
$ws.Cells.Item(5, 5).Value = -0.092952013015747
$ws.Cells.Item(5, 6).Value = 2.003179427872321
$ws.Cells.Item(5, 7).Value = -3.029003118596394
$ws.Cells.Item(5, 8).Value = -6.035639041807577

# row 6 (timestamp 400)
$ws.Cells.Item(6, 1).Value = 400
$ws.Cells.Item(6, 2).Value = "walkingToRunning"
$ws.Cells.Item(6, 3).Value = -2.974858283996582
$ws.Cells.Item(6, 4).Value = -11.37698841094971
$ws.Cells.Item(6, 5).Value = 0.3084669113159179
$ws.Cells.Item(6, 6).Value = 1.663756308759156
$ws.Cells.Item(6, 7).Value = -0.09051663155962358
$ws.Cells.Item(6, 8).Value = -2.346439961979992

# row 7 (timestamp 500)
$ws.Cells.Item(7, 1).Value = 500
$ws.Cells.Item(7, 2).Value = "walkingToRunning"
$ws.Cells.Item(7, 3).Value = -13.79527282714844
$ws.Cells.Item(7, 4).Value = -39.99485778808594
$ws.Cells.Item(7, 5).Value = 1.296123027801514
$ws.Cells.Item(7, 6).Value = -0.4272669495606185
$ws.Cells.Item(7, 7).Value = 0.4076309111423626
$ws.Cells.Item(7, 8).Value = 2.549625923720821

# row 8 (timestamp 600)
$ws.Cells.Item(8, 1).Value = 600
$ws.Cells.Item(8, 2).Value = "walkingToRunning"
$ws.Cells.Item(8, 3).Value = 0.8157472610473633
$ws.Cells.Item(8, 4).Value = -14.87919425964356
$ws.Cells.Item(8, 5).Value = 10.64366054534912
$ws.Cells.Item(8, 6).Value = -5.007111654049011
$ws.Cells.Item(8, 7).Value = 2.933262512451288
$ws.Cells.Item(8, 8).Value = 6.187744315077628

# row 9 (timestamp 700)
$ws.Cells.Item(9, 1).Value = 700
$ws.Cells.Item(9, 2).Value = "walkingToRunning"
$ws.Cells.Item(9, 3).Value = 21.28074264526367
$ws.Cells.Item(9, 4).Value = -82.66264343261719
$ws.Cells.Item(9, 5).Value = 18.65274810791016
$ws.Cells.Item(9, 6).Value = -0.8147307605277296
$ws.Cells.Item(9, 7).Value = 2.59154647297968
$ws.Cells.Item(9, 8).Value = 3.190694384458493

# row 10 (timestamp 800)
$ws.Cells.Item(10, 1).Value = 800
$ws.Cells.Item(10, 2).Value = "walkingToRunning"
$ws.Cells.Item(10, 3).Value = -67.10326385498047
$ws.Cells.Item(10, 4).Value = 36.21846008300781
$ws.Cells.Item(10, 5).Value = -4.591421127319336
$ws.Cells.Item(10, 6).Value = 4.642261807511497
$ws.Cells.Item(10, 7).Value = -2.687124653560384
$ws.Cells.Item(10, 8).Value = -4.055825431172487

# row 11 (timestamp 900)
$ws.Cells.Item(11, 1).Value = 900
$ws.Cells.Item(11, 2).Value = "walkingToRunning"
$ws.Cells.Item(11, 3).Value = -4.748800754547119
$ws.Cells.Item(11, 4).Value = -14.3919506072998
$ws.Cells.Item(11, 5).Value = 4.991169929504395
$ws.Cells.Item(11, 6).Value = 3.307725646146919
$ws.Cells.Item(11, 7).Value = -4.475619191076696
$ws.Cells.Item(11, 8).Value = -7.178602584978427

# row 12 (timestamp 1000)
$ws.Cells.Item(12, 1).Value = 1000
$ws.Cells.Item(12, 2).Value = "walkingToRunning"
$ws.Cells.Item(12, 3).Value = -13.24984931945801
$ws.Cells.Item(12, 4).Value = -18.9327278137207
$ws.Cells.Item(12, 5).Value = -1.615734577178955
$ws.Cells.Item(12, 6).Value = -0.5536176257017136
$ws.Cells.Item(12, 7).Value = -4.48860380126209
$ws.Cells.Item(12, 8).Value = -3.862921412398188

# row 13 (timestamp 1100)
$ws.Cells.Item(13, 1).Value = 1100
$ws.Cells.Item(13, 2).Value = "walkingToRunning"
$ws.Cells.Item(13, 3).Value = 31.57485389709473
$ws.Cells.Item(13, 4).Value = -14.13540840148926
$ws.Cells.Item(13, 5).Value = 37.62393188476562
$ws.Cells.Item(13, 6).Value = -7.221588652308443
$ws.Cells.Item(13, 7).Value = -2.411513553886869
$ws.Cells.Item(13, 8).Value = 2.527892092379122

# row 14 (timestamp 1200)
$ws.Cells.Item(14, 1).Value = 1200
$ws.Cells.Item(14, 2).Value = "walkingToRunning"
$ws.Cells.Item(14, 3).Value = -3.522989749908448
$ws.Cells.Item(14, 4).Value = -10.15239906311035
$ws.Cells.Item(14, 5).Value = 17.25504112243652
$ws.Cells.Item(14, 6).Value = -5.599334018986253
$ws.Cells.Item(14, 7).Value = -3.780423902883811
$ws.Cells.Item(14, 8).Value = -1.090038404232164

# row 15 (timestamp 1300)
$ws.Cells.Item(15, 1).Value = 1300
$ws.Cells.Item(15, 2).Value = "walkingToRunning"
$ws.Cells.Item(15, 3).Value = 8.178001403808594
$ws.Cells.Item(15, 4).Value = -24.92928695678711
$ws.Cells.Item(15, 5).Value = 8.257000923156738
$ws.Cells.Item(15, 6).Value = 6.376622484951461
$ws.Cells.Item(15, 7).Value = -7.650588244926631
$ws.Cells.Item(15, 8).Value = -4.948861485574298

# row 16 (timestamp 1400)
$ws.Cells.Item(16, 1).Value = 1400
$ws.Cells.Item(16, 2).Value = "walkingToRunning"
$ws.Cells.Item(16, 3).Value = -5.383133888244629
$ws.Cells.Item(16, 4).Value = 1.152879953384399
$ws.Cells.Item(16, 5).Value = 1.30066442489624
$ws.Cells.Item(16, 6).Value = 2.869635640121101
$ws.Cells.Item(16, 7).Value = -1.935632589386703
$ws.Cells.Item(16, 8).Value = -2.603090690403443

# row 17 (timestamp 1500)
$ws.Cells.Item(17, 1).Value = 1500
$ws.Cells.Item(17, 2).Value = "walkingToRunning"
$ws.Cells.Item(17, 3).Value = -21.55855941772461
$ws.Cells.Item(17, 4).Value = -18.73270034790039
$ws.Cells.Item(17, 5).Value = -9.058673858642578
$ws.Cells.Item(17, 6).Value = 2.96556776325869
$ws.Cells.Item(17, 7).Value = 3.462880915984825
$ws.Cells.Item(17, 8).Value = -0.1149900999010921

# row 18 (timestamp 1600)
$ws.Cells.Item(18, 1).Value = 1600
$ws.Cells.Item(18, 2).Value = "walkingToRunning"
$ws.Cells.Item(18, 3).Value = 4.207836627960205
$ws.Cells.Item(18, 4).Value = -60.38365173339844
$ws.Cells.Item(18, 5).Value = 21.30802917480469
$ws.Cells.Item(18, 6).Value = -4.452823406312537
$ws.Cells.Item(18, 7).Value = 0.1738118369405208
$ws.Cells.Item(18, 8).Value = 2.291137044022761

# row 19 (timestamp 1700)
$ws.Cells.Item(19, 1).Value = 1700
$ws.Cells.Item(19, 2).Value = "walkingToRunning"
$ws.Cells.Item(19, 3).Value = -7.397396087646484
$ws.Cells.Item(19, 4).Value = -7.57304859161377
$ws.Cells.Item(19, 5).Value = 8.563810348510742
$ws.Cells.Item(19, 6).Value = -4.949691057205192
$ws.Cells.Item(19, 7).Value = 5.274479982329643
$ws.Cells.Item(19, 8).Value = 2.80145074390779

# row 20 (timestamp 1800)
$ws.Cells.Item(20, 1).Value = 1800
$ws.Cells.Item(20, 2).Value = "walkingToRunning"
$ws.Cells.Item(20, 3).Value = 16.9281063079834
$ws.Cells.Item(20, 4).Value = -77.95144653320312
$ws.Cells.Item(20, 5).Value = 60.47911071777344
$ws.Cells.Item(20, 6).Value = -2.519368160061688
$ws.Cells.Item(20, 7).Value = 9.493784741657493
$ws.Cells.Item(20, 8).Value = -1.797140115644872

# row 21 (timestamp 1900)
$ws.Cells.Item(21, 1).Value = 1900
$ws.Cells.Item(21, 2).Value = "walkingToRunning"
$ws.Cells.Item(21, 3).Value = -18.38624954223633
$ws.Cells.Item(21, 4).Value = 4.694716930389404
$ws.Cells.Item(21, 5).Value = -15.62700939178467
$ws.Cells.Item(21, 6).Value = 1.409833646402128
$ws.Cells.Item(21, 7).Value = 1.814598339359917
$ws.Cells.Item(21, 8).Value = -1.559189867682572

# row 22 (timestamp 2000)
$ws.Cells.Item(22, 1).Value = 2000
$ws.Cells.Item(22, 2).Value = "walkingToRunning"
$ws.Cells.Item(22, 3).Value = 36.83388900756836
$ws.Cells.Item(22, 4).Value = -10.96193885803223
$ws.Cells.Item(22, 5).Value = 0.4498906135559082
$ws.Cells.Item(22, 6).Value = 2.461990158732359
$ws.Cells.Item(22, 7).Value = -6.157769249706776
$ws.Cells.Item(22, 8).Value = -2.199485290341278

# row 23 (timestamp 2100)
$ws.Cells.Item(23, 1).Value = 2100
$ws.Cells.Item(23, 2).Value = "walkingToRunning"
$ws.Cells.Item(23, 3).Value = -8.928971290588379
$ws.Cells.Item(23, 4).Value = -17.86810111999512
$ws.Cells.Item(23, 5).Value = 8.281005859375
$ws.Cells.Item(23, 6).Value = 2.753937654378926
$ws.Cells.Item(23, 7).Value = -9.407423193861781
$ws.Cells.Item(23, 8).Value = -3.08937735819234

# row 24 (timestamp 2200)
$ws.Cells.Item(24, 1).Value = 2200
$ws.Cells.Item(24, 2).Value = "walkingToRunning"
$ws.Cells.Item(24, 3).Value = -14.49608421325684
$ws.Cells.Item(24, 4).Value = -1.527808666229248
$ws.Cells.Item(24, 5).Value = 44.4189453125
$ws.Cells.Item(24, 6).Value = -4.065274791019721
$ws.Cells.Item(24, 7).Value = -1.80112353185328
$ws.Cells.Item(24, 8).Value = 2.134828872796967

# row 25 (timestamp 2300)
$ws.Cells.Item(25, 1).Value = 2300
$ws.Cells.Item(25, 2).Value = "walkingToRunning"
$ws.Cells.Item(25, 3).Value = -12.06443023681641
$ws.Cells.Item(25, 4).Value = 6.844409942626953
$ws.Cells.Item(25, 5).Value = 19.9449577331543
$ws.Cells.Item(25, 6).Value = -7.774596919373714
$ws.Cells.Item(25, 7).Value = 2.037927262666848
$ws.Cells.Item(25, 8).Value = 2.593074496199395

# row 26 (timestamp 2400)
$ws.Cells.Item(26, 1).Value = 2400
$ws.Cells.Item(26, 2).Value = "walkingToRunning"
$ws.Cells.Item(26, 3).Value = 6.954762935638428
$ws.Cells.Item(26, 4).Value = -76.15243530273438
$ws.Cells.Item(26, 5).Value = 24.18494606018066
$ws.Cells.Item(26, 6).Value = 0.1093276535592701
$ws.Cells.Item(26, 7).Value = 10.26185343905173
$ws.Cells.Item(26, 8).Value = -3.281300154764442

# row 27 (timestamp 2500)
$ws.Cells.Item(27, 1).Value = 2500
$ws.Cells.Item(27, 2).Value = "walkingToRunning"
$ws.Cells.Item(27, 3).Value = 6.384909629821777
$ws.Cells.Item(27, 4).Value = 5.00542688369751
$ws.Cells.Item(27, 5).Value = -29.23712921142578
$ws.Cells.Item(27, 6).Value = 3.391102220953991
$ws.Cells.Item(27, 7).Value = -5.063158106513092
$ws.Cells.Item(27, 8).Value = -1.229867340406294

# row 28 (timestamp 2600)
$ws.Cells.Item(28, 1).Value = 2600
$ws.Cells.Item(28, 2).Value = "walkingToRunning"
$ws.Cells.Item(28, 3).Value = -34.79932403564453
$ws.Cells.Item(28, 4).Value = -7.816071510314941
$ws.Cells.Item(28, 5).Value = 1.089200496673584
$ws.Cells.Item(28, 6).Value = 6.242717754550078
$ws.Cells.Item(28, 7).Value = -0.2678701499613687
$ws.Cells.Item(28, 8).Value = -4.157948156682439

# row 29 (timestamp 2700)
$ws.Cells.Item(29, 1).Value = 2700
$ws.Cells.Item(29, 2).Value = "walkingToRunning"
$ws.Cells.Item(29, 3).Value = -17.0820198059082
$ws.Cells.Item(29, 4).Value = -31.8654670715332
$ws.Cells.Item(29, 5).Value = 12.90904235839844
$ws.Cells.Item(29, 6).Value = 1.95433324720798
$ws.Cells.Item(29, 7).Value = -6.421389656095905
$ws.Cells.Item(29, 8).Value = 0.3005734565781801

# row 30 (timestamp 2800)
$ws.Cells.Item(30, 1).Value = 2800
$ws.Cells.Item(30, 2).Value = "walkingToRunning"
$ws.Cells.Item(30, 3).Value = 2.159783363342285
$ws.Cells.Item(30, 4).Value = 0.4922776222229004
$ws.Cells.Item(30, 5).Value = 7.778494358062744
$ws.Cells.Item(30, 6).Value = -2.448658175584859
$ws.Cells.Item(30, 7).Value = -1.071053583447552
$ws.Cells.Item(30, 8).Value = 4.803342400527646

# row 31 (timestamp 2900)
$ws.Cells.Item(31, 1).Value = 2900
$ws.Cells.Item(31, 2).Value = "walkingToRunning"
$ws.Cells.Item(31, 3).Value = 3.85674524307251
$ws.Cells.Item(31, 4).Value = 1.991205930709839
$ws.Cells.Item(31, 5).Value = 21.4826774597168
$ws.Cells.Item(31, 6).Value = -7.789634487977738
$ws.Cells.Item(31, 7).Value = -0.2244347770038115
$ws.Cells.Item(31, 8).Value = 7.027555852401433
